# book_cover.pptx — "finished half of chpt 6"
#
# 1) Refresh the cached "datetimeFigureOut" date field (11/3/17 -> 11/17/17)
#    on the slide master and every slide layout's Date placeholder.
# 2) Split the subtitle run "         DOUGLAS RUBIN PhD" so the extra
#    leading space is rendered smaller (sz 2000) before the name.
# 3) Nudge the three lower textboxes 7pt to the right; narrow the title
#    textbox to make room, and prefix the heading with "A ".

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "11/3/17") {
                $tr.Text = "11/17/17"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholders $layouts.Item($i)
}

# --- Slide 1 ---------------------------------------------------------
$s = $p.Slides.Item(1)

# Title 1: "         DOUGLAS RUBIN PhD" -> extra small space + name
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "          DOUGLAS RUBIN PhD"
$extraSpace = $titleRange.Characters(10, 1)
$extraSpace.Font.Size = 20

# TextBox 5: heading textbox - move/resize + prefix text with "A "
$heading = $s.Shapes.Item(3)
$headingFirstRun = $heading.TextFrame.TextRange.Characters(1, 28)
$headingFirstRun.Text = "A Complete Solutions Guide to "
$heading.Left = 75.2696876525879
$heading.Width = 475.0636138916016

# TextBox 6: subtitle textbox - move only
$subtitle = $s.Shapes.Item(4)
$subtitle.Left = 75.26976394653322

# TextBox 7: accent bar - move only
$bar = $s.Shapes.Item(5)
$bar.Left = 57.269685745239265
